# Update "想去人数" (want-to-go count) values in column F on the sheets
# that list individual events: "展览" (sheet1) and "全部类型" (sheet4).
# Sheets "演出" and "本地生活" only contain header rows, so nothing to change there.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    3  = 2142
    4  = 1629
    5  = 321
    6  = 1040
    7  = 542
    9  = 5682
    10 = 80
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
